$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Delete Group"

# Update header text
$ws.Range("A1").Value = "Delete Name"

# Remove column B entirely (was "Item Group Name")
$ws.Columns.Item(2).Delete()

# Update selection to match target (whole column A selected)
$ws.Range("A1:A1048576").Select()
